$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Customer Records" to "Employee Records"
$ws.Name = "Employee Records"

# Expand the table with three new columns: Email, Salary, Birthday
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Add() | Out-Null
$tbl.ListColumns.Add() | Out-Null
$tbl.ListColumns.Add() | Out-Null

$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Salary"
$ws.Range("F1").Value = "Birthday"

# Force the Birthday column to store plain text so "2020-12-12" is not
# reinterpreted as a date serial number.
$ws.Range("F2:F10").NumberFormat = "@"

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = "andres@gmail.com"
    $ws.Cells.Item($r, 5).Value = 1500000
    $ws.Cells.Item($r, 6).Value = "2020-12-12"
}
